$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$districtNames = @(
    "Kalangala",
    "Kampala",
    "Kiboga",
    "Luwero",
    "Masaka",
    "Mpigi",
    "Mubende",
    "Mukono",
    "Nakasongola",
    "Rakai",
    "Ssembabule",
    "Kayunga",
    "Wakiso",
    "Lyantonde",
    "Mityana",
    "Nakaseke",
    "Buikwe",
    "Bukomansimbi",
    "Butambala",
    "Buvuma",
    "Gomba",
    "Kalungu",
    "Kyankwanzi",
    "Lwengo",
    "Bugiri",
    "Busia",
    "Iganga",
    "Jinja",
    "Kamuli",
    "Kapchorwa",
    "Katakwi",
    "Kumi",
    "Mbale",
    "Pallisa",
    "Soroti",
    "Tororo",
    "Kaberamaido",
    "Mayuge",
    "Sironko",
    "Amuria",
    "Budaka",
    "Bududa",
    "Bukedea",
    "Bukwo",
    "Butaleja",
    "Kaliro",
    "Manafwa",
    "Namutumba",
    "Bulambuli",
    "Buyende",
    "Kibuku",
    "Kween",
    "Luuka",
    "Namayingo",
    "Ngora",
    "Serere",
    "Adjumani",
    "Apac",
    "Arua",
    "Gulu",
    "Kitgum",
    "Kotido",
    "Lira",
    "Moroto",
    "Moyo",
    "Nebbi",
    "Nakapiripirit",
    "Pader",
    "Yumbe",
    "Abim",
    "Amolatar",
    "Amuru",
    "Dokolo",
    "Kaabong",
    "Koboko",
    "Maracha",
    "Oyam",
    "Agago",
    "Alebtong",
    "Amudat",
    "Kole",
    "Lamwo",
    "Napak",
    "Nwoya",
    "Otuke",
    "Zombo",
    "Bundibugyo",
    "Bushenyi",
    "Hoima",
    "Kabale",
    "Kabarole",
    "Kasese",
    "Kibaale",
    "Kisoro",
    "Masindi",
    "Mbarara",
    "Ntungamo",
    "Rukungiri",
    "Kamwenge",
    "Kanungu",
    "Kyenjojo",
    "Buliisa",
    "Ibanda",
    "Isingiro",
    "Kiruhura",
    "Buhweju",
    "Kiryandongo",
    "Kyegegwa",
    "Mitooma",
    "Ntoroko",
    "Rubirizi",
    "Sheema"
)

for ($i = 0; $i -lt $districtNames.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $districtNames[$i]
}